# Update numeric values (imputed results) in the active worksheet.
# These correspond to updated RandomForest imputation output values
# for columns B and D across several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.129199999999996
$ws.Range("B6").Value = 5.405000000000005
$ws.Range("B7").Value = 5.321100000000004
$ws.Range("D7").Value = -7.328000000000003
$ws.Range("D12").Value = -7.211699999999999
$ws.Range("D15").Value = -8.858199999999993
$ws.Range("B16").Value = 6.7205
$ws.Range("B20").Value = 9.451799999999993
$ws.Range("D20").Value = -7.973799999999992
$ws.Range("D21").Value = -7.921899999999991
$ws.Range("D22").Value = -7.560300000000001
$ws.Range("D23").Value = -7.214199999999997
$ws.Range("B28").Value = 6.339600000000001
$ws.Range("B29").Value = 5.067800000000004
$ws.Range("D29").Value = -7.318599999999998
$ws.Range("B32").Value = 7.530799999999995
$ws.Range("D34").Value = -7.749700000000003
$ws.Range("B40").Value = 8.964499999999996
$ws.Range("D42").Value = -8.045700000000004
$ws.Range("D43").Value = -8.279300000000001
$ws.Range("D44").Value = -7.536400000000002
$ws.Range("D45").Value = -7.806999999999999
$ws.Range("B46").Value = 6.123399999999998
$ws.Range("D46").Value = -8.238200000000001
$ws.Range("D50").Value = -8.106999999999996
$ws.Range("B51").Value = 5.0429
$ws.Range("D51").Value = -7.465699999999999
$ws.Range("B52").Value = 5.060700000000001
$ws.Range("B57").Value = 5.487199999999994
$ws.Range("B59").Value = 4.657799999999997
$ws.Range("B62").Value = 5.6377
$ws.Range("B66").Value = 5.386599999999998
$ws.Range("D66").Value = -7.347800000000003
$ws.Range("D67").Value = -6.4654
$ws.Range("B73").Value = 8.331600000000002
$ws.Range("B74").Value = 9.349599999999992
$ws.Range("D79").Value = -6.358300000000003
$ws.Range("D84").Value = -8.858299999999996
$ws.Range("B92").Value = 4.835299999999998
$ws.Range("D92").Value = -6.411600000000001
$ws.Range("D97").Value = -8.742900000000002
$ws.Range("B100").Value = 5.5641
